$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the formatting of the existing
# header cells (bold, centered, bordered) so it reuses the same style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column's data values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
